$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Submit orders")
$ws.Cells.Item(102, 1).Value = "11.28.2022 12:03 (Kyiv+Israel) 10:03 (UTC) 19:03 (Japan) 15:33 (India)"
$ws.Cells.Item(102, 2).Value = 0.973
$ws.Cells.Item(102, 3).Value = -0.246
$ws.Cells.Item(102, 4).Value = "***"
$ws.Cells.Item(102, 5).Value = "***"
$ws.Cells.Item(103, 1).Value = "12.05.2022 09:21 (Kyiv+Israel) 07:21 (UTC) 16:21 (Japan) 12:51 (India)"
$ws.Cells.Item(103, 2).Value = 2.624
$ws.Cells.Item(103, 3).Value = -1.897
$ws.Cells.Item(103, 4).Value = "***"
$ws.Cells.Item(103, 5).Value = "***"
$ws.Cells.Item(104, 1).Value = "12.05.2022 10:22 (Kyiv+Israel) 08:22 (UTC) 17:22 (Japan) 13:52 (India)"
$ws.Cells.Item(104, 2).Value = 1.808
$ws.Cells.Item(104, 3).Value = -1.081
$ws.Cells.Item(104, 4).Value = "***"
$ws.Cells.Item(104, 5).Value = "***"
$ws.Cells.Item(105, 1).Value = "12.05.2022 11:29 (Kyiv+Israel) 09:29 (UTC) 18:29 (Japan) 14:59 (India)"
$ws.Cells.Item(105, 2).Value = 3.131
$ws.Cells.Item(105, 3).Value = -2.404
$ws.Cells.Item(105, 4).Value = "***"
$ws.Cells.Item(105, 5).Value = "***"
$ws.Cells.Item(106, 1).Value = "12.05.2022 17:14 (Kyiv+Israel) 15:14 (UTC) 00:14 (Japan) 20:44 (India)"
$ws.Cells.Item(106, 2).Value = "***"
$ws.Cells.Item(106, 3).Value = "***"
$ws.Cells.Item(106, 4).Value = 1.843
$ws.Cells.Item(106, 5).Value = -0.7050000000000001
$ws.Cells.Item(107, 1).Value = "12.05.2022 17:15 (Kyiv+Israel) 15:15 (UTC) 00:15 (Japan) 20:45 (India)"
$ws.Cells.Item(107, 2).Value = 2.408
$ws.Cells.Item(107, 3).Value = -1.681
$ws.Cells.Item(107, 4).Value = "***"
$ws.Cells.Item(107, 5).Value = "***"

$ws = $wb.Worksheets.Item("Submit internet survey")
$ws.Cells.Item(100, 1).Value = "12.05.2022 09:09 (Kyiv+Israel) 07:09 (UTC) 16:09 (Japan) 12:39 (India)"
$ws.Cells.Item(100, 2).Value = "***"
$ws.Cells.Item(100, 3).Value = "***"
$ws.Cells.Item(100, 4).Value = 0.998
$ws.Cells.Item(100, 5).Value = -0.405
$ws.Cells.Item(101, 1).Value = "12.05.2022 09:24 (Kyiv+Israel) 07:24 (UTC) 16:24 (Japan) 12:54 (India)"
$ws.Cells.Item(101, 2).Value = 0.983
$ws.Cells.Item(101, 3).Value = -0.403
$ws.Cells.Item(101, 4).Value = "***"
$ws.Cells.Item(101, 5).Value = "***"
$ws.Cells.Item(102, 1).Value = "12.05.2022 10:24 (Kyiv+Israel) 08:24 (UTC) 17:24 (Japan) 13:54 (India)"
$ws.Cells.Item(102, 2).Value = 0.9389999999999999
$ws.Cells.Item(102, 3).Value = -0.359
$ws.Cells.Item(102, 4).Value = "***"
$ws.Cells.Item(102, 5).Value = "***"
$ws.Cells.Item(103, 1).Value = "12.05.2022 11:32 (Kyiv+Israel) 09:32 (UTC) 18:32 (Japan) 15:02 (India)"
$ws.Cells.Item(103, 2).Value = 1.162
$ws.Cells.Item(103, 3).Value = -0.582
$ws.Cells.Item(103, 4).Value = "***"
$ws.Cells.Item(103, 5).Value = "***"
$ws.Cells.Item(104, 1).Value = "12.05.2022 17:18 (Kyiv+Israel) 15:18 (UTC) 00:18 (Japan) 20:48 (India)"
$ws.Cells.Item(104, 2).Value = "***"
$ws.Cells.Item(104, 3).Value = "***"
$ws.Cells.Item(104, 4).Value = 1.371
$ws.Cells.Item(104, 5).Value = -0.778
$ws.Cells.Item(105, 1).Value = "12.05.2022 17:19 (Kyiv+Israel) 15:19 (UTC) 00:19 (Japan) 20:49 (India)"
$ws.Cells.Item(105, 2).Value = 1.408
$ws.Cells.Item(105, 3).Value = -0.828
$ws.Cells.Item(105, 4).Value = "***"
$ws.Cells.Item(105, 5).Value = "***"

$ws = $wb.Worksheets.Item("Submit a phone survey")
$ws.Cells.Item(92, 1).Value = "12.05.2022 10:09 (Kyiv+Israel) 08:09 (UTC) 17:09 (Japan) 13:39 (India)"
$ws.Cells.Item(92, 2).Value = 3.086
$ws.Cells.Item(92, 3).Value = -1.982
$ws.Cells.Item(92, 4).Value = "***"
$ws.Cells.Item(92, 5).Value = "***"
$ws.Cells.Item(93, 1).Value = "12.05.2022 10:27 (Kyiv+Israel) 08:27 (UTC) 17:27 (Japan) 13:57 (India)"
$ws.Cells.Item(93, 2).Value = 2.8
$ws.Cells.Item(93, 3).Value = -1.696
$ws.Cells.Item(93, 4).Value = "***"
$ws.Cells.Item(93, 5).Value = "***"
$ws.Cells.Item(94, 1).Value = "12.05.2022 17:22 (Kyiv+Israel) 15:22 (UTC) 00:22 (Japan) 20:52 (India)"
$ws.Cells.Item(94, 2).Value = 8.101000000000001
$ws.Cells.Item(94, 3).Value = -6.997000000000001
$ws.Cells.Item(94, 4).Value = "***"
$ws.Cells.Item(94, 5).Value = "***"
$ws.Cells.Item(95, 1).Value = "12.05.2022 23:30 (Kyiv+Israel) 21:30 (UTC) 06:30 (Japan) 03:00 (India)"
$ws.Cells.Item(95, 2).Value = 2.945
$ws.Cells.Item(95, 3).Value = -1.841
$ws.Cells.Item(95, 4).Value = "***"
$ws.Cells.Item(95, 5).Value = "***"
$ws.Cells.Item(96, 1).Value = "12.05.2022 23:40 (Kyiv+Israel) 21:40 (UTC) 06:40 (Japan) 03:10 (India)"
$ws.Cells.Item(96, 2).Value = "***"
$ws.Cells.Item(96, 3).Value = "***"
$ws.Cells.Item(96, 4).Value = 1.833
$ws.Cells.Item(96, 5).Value = -0.2489999999999999
$ws.Cells.Item(97, 1).Value = "12.05.2022 23:46 (Kyiv+Israel) 21:46 (UTC) 06:46 (Japan) 03:16 (India)"
$ws.Cells.Item(97, 2).Value = "***"
$ws.Cells.Item(97, 3).Value = "***"
$ws.Cells.Item(97, 4).Value = 1.666
$ws.Cells.Item(97, 5).Value = -0.08199999999999985

$ws = $wb.Worksheets.Item("Checkertificate")
$ws.Cells.Item(104, 1).Value = "12.05.2022 09:53 (Kyiv+Israel) 07:53 (UTC) 16:53 (Japan) 13:23 (India)"
$ws.Cells.Item(104, 2).Value = 0.831
$ws.Cells.Item(104, 3).Value = -0.1659999999999999
$ws.Cells.Item(104, 4).Value = "***"
$ws.Cells.Item(104, 5).Value = "***"
$ws.Cells.Item(105, 1).Value = "12.05.2022 09:55 (Kyiv+Israel) 07:55 (UTC) 16:55 (Japan) 13:25 (India)"
$ws.Cells.Item(105, 2).Value = "***"
$ws.Cells.Item(105, 3).Value = "***"
$ws.Cells.Item(105, 4).Value = 1.211
$ws.Cells.Item(105, 5).Value = -0.289
$ws.Cells.Item(106, 1).Value = "12.05.2022 10:17 (Kyiv+Israel) 08:17 (UTC) 17:17 (Japan) 13:47 (India)"
$ws.Cells.Item(106, 2).Value = 1.025
$ws.Cells.Item(106, 3).Value = -0.3599999999999999
$ws.Cells.Item(106, 4).Value = "***"
$ws.Cells.Item(106, 5).Value = "***"
$ws.Cells.Item(107, 1).Value = "12.05.2022 11:39 (Kyiv+Israel) 09:39 (UTC) 18:39 (Japan) 15:09 (India)"
$ws.Cells.Item(107, 2).Value = 0.974
$ws.Cells.Item(107, 3).Value = -0.3089999999999999
$ws.Cells.Item(107, 4).Value = "***"
$ws.Cells.Item(107, 5).Value = "***"
$ws.Cells.Item(108, 1).Value = "12.05.2022 17:24 (Kyiv+Israel) 15:24 (UTC) 00:24 (Japan) 20:54 (India)"
$ws.Cells.Item(108, 2).Value = "***"
$ws.Cells.Item(108, 3).Value = "***"
$ws.Cells.Item(108, 4).Value = 1.233
$ws.Cells.Item(108, 5).Value = -0.3110000000000001
$ws.Cells.Item(109, 1).Value = "12.05.2022 17:32 (Kyiv+Israel) 15:32 (UTC) 00:32 (Japan) 21:02 (India)"
$ws.Cells.Item(109, 2).Value = 0.889
$ws.Cells.Item(109, 3).Value = -0.224
$ws.Cells.Item(109, 4).Value = "***"
$ws.Cells.Item(109, 5).Value = "***"
